$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 177.6923
$ws.Cells.Item(2, 9).Value = 222.5
$ws.Cells.Item(2, 10).Value = 106
$ws.Cells.Item(2, 11).Value = 222.5
$ws.Cells.Item(2, 12).Value = 106
$ws.Cells.Item(2, 13).Value = -109.5
$ws.Cells.Item(2, 14).Value = -332
# Row 6
$ws.Cells.Item(6, 8).Value = 335.9091
$ws.Cells.Item(6, 9).Value = 294.5
$ws.Cells.Item(6, 10).Value = 750
$ws.Cells.Item(6, 11).Value = 883.5
$ws.Cells.Item(6, 12).Value = 2250
$ws.Cells.Item(6, 13).Value = -771.5
$ws.Cells.Item(6, 14).Value = -2474
# Row 8
$ws.Cells.Item(8, 8).Value = 3279.25
$ws.Cells.Item(8, 9).Value = 3431.2
$ws.Cells.Item(8, 10).Value = 1000
$ws.Cells.Item(8, 11).Value = 10293.6
$ws.Cells.Item(8, 12).Value = 3000
$ws.Cells.Item(8, 13).Value = -10154.6
$ws.Cells.Item(8, 14).Value = -3278
# Row 95
$ws.Cells.Item(95, 8).Value = 80000
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 80000
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 80000
$ws.Cells.Item(95, 14).Value = -85492

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 14
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 14).ClearContents()
# Row 61
$ws.Cells.Item(61, 8).Value = 1598.6923
$ws.Cells.Item(61, 9).Value = 1092.8889
$ws.Cells.Item(61, 10).Value = 2736.75
$ws.Cells.Item(61, 11).Value = 1092.8889
$ws.Cells.Item(61, 12).Value = 2736.75
$ws.Cells.Item(61, 13).Value = -880.8888999999999
$ws.Cells.Item(61, 14).Value = -3160.75
# Row 74
$ws.Cells.Item(74, 8).Value = 914.17645
$ws.Cells.Item(74, 9).Value = 914.17645
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 914.17645
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -40.17645000000005
# Row 77
$ws.Cells.Item(77, 8).Value = 914.17645
$ws.Cells.Item(77, 9).Value = 914.17645
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 4570.882250000001
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -202.8822500000006
# Row 134
$ws.Cells.Item(134, 8).Value = 38000
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 38000
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 38000
$ws.Cells.Item(134, 14).Value = -48140
# Row 136
$ws.Cells.Item(136, 8).Value = 1598.6923
$ws.Cells.Item(136, 9).Value = 1092.8889
$ws.Cells.Item(136, 10).Value = 2736.75
$ws.Cells.Item(136, 11).Value = 3278.6667
$ws.Cells.Item(136, 12).Value = 8210.25
$ws.Cells.Item(136, 13).Value = -728.6666999999998
$ws.Cells.Item(136, 14).Value = -13310.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Cells.Item(7, 8).Value = 2200
$ws.Cells.Item(7, 9).Value = 600
$ws.Cells.Item(7, 10).Value = 3000
$ws.Cells.Item(7, 11).Value = 600
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = -487
$ws.Cells.Item(7, 14).Value = -3226
# Row 16
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Cells.Item(12, 8).Value = 3006
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 3006
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 3006
$ws.Cells.Item(12, 14).Value = -3346
$ws.Cells.Item(12, 13).ClearContents()
# Row 31
$ws.Cells.Item(31, 8).Value = 2896.6
$ws.Cells.Item(31, 9).Value = 1751.619
$ws.Cells.Item(31, 10).Value = 5568.222
$ws.Cells.Item(31, 11).Value = 1751.619
$ws.Cells.Item(31, 12).Value = 5568.222
$ws.Cells.Item(31, 13).Value = -1456.619
$ws.Cells.Item(31, 14).Value = -6158.222
# Row 34
$ws.Cells.Item(34, 8).Value = 2896.6
$ws.Cells.Item(34, 9).Value = 1751.619
$ws.Cells.Item(34, 10).Value = 5568.222
$ws.Cells.Item(34, 11).Value = 1751.619
$ws.Cells.Item(34, 12).Value = 5568.222
$ws.Cells.Item(34, 13).Value = -1549.619
$ws.Cells.Item(34, 14).Value = -5972.222
# Row 95
$ws.Cells.Item(95, 8).Value = 12592
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 12592
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 12592
$ws.Cells.Item(95, 14).Value = -18084
# Row 105
$ws.Cells.Item(105, 8).Value = 4570
$ws.Cells.Item(105, 9).Value = 4781.4287
$ws.Cells.Item(105, 10).Value = 4200
$ws.Cells.Item(105, 11).Value = 4781.4287
$ws.Cells.Item(105, 12).Value = 4200
$ws.Cells.Item(105, 13).Value = -3034.4287
$ws.Cells.Item(105, 14).Value = -7694
# Row 132
$ws.Cells.Item(132, 8).Value = 1609.5555
$ws.Cells.Item(132, 9).Value = 1185.875
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 3557.625
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -1027.625
$ws.Cells.Item(132, 14).Value = -20057
# Row 134
$ws.Cells.Item(134, 8).Value = 1316.3334
$ws.Cells.Item(134, 9).Value = 885.3043
$ws.Cells.Item(134, 10).Value = 2732.5715
$ws.Cells.Item(134, 11).Value = 2655.9129
$ws.Cells.Item(134, 12).Value = 8197.7145
$ws.Cells.Item(134, 13).Value = -120.9129000000003
$ws.Cells.Item(134, 14).Value = -13267.7145

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 7913.923
$ws.Cells.Item(7, 9).Value = 9218.272000000001
$ws.Cells.Item(7, 10).Value = 740
$ws.Cells.Item(7, 11).Value = 27654.816
$ws.Cells.Item(7, 12).Value = 2220
$ws.Cells.Item(7, 13).Value = -27542.816
$ws.Cells.Item(7, 14).Value = -2444
# Row 19
$ws.Cells.Item(19, 8).Value = 996.6667
$ws.Cells.Item(19, 9).Value = 980
# Row 131
$ws.Cells.Item(131, 8).Value = 754.4167
$ws.Cells.Item(131, 9).Value = 259.2857
$ws.Cells.Item(131, 10).Value = 958.2941
$ws.Cells.Item(131, 11).Value = 777.8571000000001
$ws.Cells.Item(131, 12).Value = 2874.8823
$ws.Cells.Item(131, 13).Value = 4262.1429
$ws.Cells.Item(131, 14).Value = -12954.8823

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Cells.Item(31, 8).Value = 3074.2727
$ws.Cells.Item(31, 9).Value = 778.2
$ws.Cells.Item(31, 10).Value = 26035
$ws.Cells.Item(31, 11).Value = 778.2
$ws.Cells.Item(31, 12).Value = 26035
$ws.Cells.Item(31, 13).Value = -486.2
$ws.Cells.Item(31, 14).Value = -26619
# Row 37
$ws.Cells.Item(37, 8).Value = 3074.2727
$ws.Cells.Item(37, 9).Value = 778.2
$ws.Cells.Item(37, 10).Value = 26035
$ws.Cells.Item(37, 11).Value = 778.2
$ws.Cells.Item(37, 12).Value = 26035
$ws.Cells.Item(37, 13).Value = -501.2
$ws.Cells.Item(37, 14).Value = -26589
# Row 132
$ws.Cells.Item(132, 8).Value = 4543.4
$ws.Cells.Item(132, 9).Value = 4802.0356
$ws.Cells.Item(132, 10).Value = 3508.8572
$ws.Cells.Item(132, 11).Value = 14406.1068
$ws.Cells.Item(132, 12).Value = 10526.5716
$ws.Cells.Item(132, 13).Value = -11876.1068
$ws.Cells.Item(132, 14).Value = -15586.5716

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Cells.Item(9, 8).Value = 348
$ws.Cells.Item(9, 9).Value = 348
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 348
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = -124
$ws.Cells.Item(9, 14).ClearContents()
# Row 16
$ws.Cells.Item(16, 8).Value = 1026.4667
$ws.Cells.Item(16, 9).Value = 1178.7778
$ws.Cells.Item(16, 10).Value = 798
$ws.Cells.Item(16, 11).Value = 1178.7778
$ws.Cells.Item(16, 12).Value = 798
$ws.Cells.Item(16, 13).Value = -1008.7778
$ws.Cells.Item(16, 14).Value = -1138
# Row 30
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).ClearContents()
$ws.Cells.Item(30, 14).ClearContents()
# Row 40
$ws.Cells.Item(40, 8).Value = 2194.889
$ws.Cells.Item(40, 9).Value = 2100.5334
$ws.Cells.Item(40, 10).Value = 2666.6667
$ws.Cells.Item(40, 11).Value = 2100.5334
$ws.Cells.Item(40, 12).Value = 2666.6667
$ws.Cells.Item(40, 13).Value = -1964.5334
$ws.Cells.Item(40, 14).Value = -2938.6667
# Row 132
$ws.Cells.Item(132, 8).Value = 3396
$ws.Cells.Item(132, 9).Value = 2450.3
$ws.Cells.Item(132, 10).Value = 4184.0835
$ws.Cells.Item(132, 11).Value = 7350.900000000001
$ws.Cells.Item(132, 12).Value = 12552.2505
$ws.Cells.Item(132, 13).Value = -4820.900000000001
$ws.Cells.Item(132, 14).Value = -17612.2505

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Cells.Item(15, 8).Value = 1000
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 1000
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 1000
$ws.Cells.Item(15, 14).Value = -1576
# Row 18
$ws.Cells.Item(18, 8).Value = 9475
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 9475
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 9475
$ws.Cells.Item(18, 14).Value = -9821
